$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto data (prices/volumes/coin swaps) per commit diff.
# Force text number format on every touched cell first so Excel does not
# auto-coerce text-looking-like-numbers (e.g. "1.00", "573.10") into floats,
# matching the original inlineStr (text) cell type in the workbook.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.814.73'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.96%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.559.50'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.30%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.01'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '188.79'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.624'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.69%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.548.80'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.32%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +15.02%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.48'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.51'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.125.47'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '70.817.76'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.15%  '
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.609.92'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.72%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.19'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.76'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '573.10'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.63%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.61'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -6.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.58'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.91'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '94.26'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.20'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.32%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.29'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.69'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.21'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.30'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.87%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.66'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.80'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +22.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.32'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +7.23%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '530.82'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '38.43'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0802'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.618.79'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +9.64%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.43'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0469'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.48'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.74%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.38%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.28'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.66%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.138'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.998'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +6.30%  '
